# -----------------------------------------------------------------------
# Extracts 2 tables and add debugging
# - Renames the existing sheet to "Member Co-Pay"
# - Fixes a couple of mis-typed cells on that sheet (E2 "nan" -> "Null",
#   I2/J2 stored as real number/boolean instead of text)
# - Adds a second sheet "Insured Co-Pay" with the same header row, but
#   with the E1 header relabelled "Insured Co-Pay"
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: rename "Sheet1" -> "Member Co-Pay" ---
$ws1.Name = "Member Co-Pay"

# --- Data-quality fixes on Sheet1 row 2 ---
# E2: literal "nan" string -> "Null"
$ws1.Range("E2").Value = "Null"

# I2: was stored as text "100" -> store as a real number
$ws1.Range("I2").Value = 100

# J2: was stored as text "TRUE" -> store as a real boolean
$ws1.Range("J2").Value = $true

# --- Add the second sheet right after "Member Co-Pay" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Insured Co-Pay"

# Match the header row's style (bold, centered, top-aligned, thin box border)
$hdrRange = $ws2.Range("A1:J1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108   # xlCenter
$hdrRange.VerticalAlignment = -4160     # xlTop
$hdrRange.Borders.LineStyle = 1         # xlContinuous
$hdrRange.Borders.Weight = 2            # xlThin

# Preserve the outline defaults (matches the source sheet's sheetPr)
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1

# Match the workbook's page margins (0.75/0.75/1/1/0.5/0.5 in, i.e. 54/54/72/72/36/36 pt)
$ws2.PageSetup.LeftMargin   = 54
$ws2.PageSetup.RightMargin  = 54
$ws2.PageSetup.TopMargin    = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$headers = @("Policy No", "Card No", "No", "Coverage", "Insured Co-Pay", "Limit", "Balance", "Coverage (raw)", "MatchScore", "Matched")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Keep "Member Co-Pay" as the active/selected sheet (matches activeTab="0")
$ws1.Activate()
$ws1.Select()
